# Auto-generated edit script: refresh cryptos price/volume data (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.179.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "'2.076.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'254.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("D7").Value = "'60.72"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.65%  "

$ws.Range("E9").Value = "  +4.87%  "

$ws.Range("D10").Value = "'61.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("E11").Value = "  +8.24%  "

$ws.Range("E12").Value = "  +2.46%  "

$ws.Range("D13").Value = "'16.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.71%  "

$ws.Range("D14").Value = "'2.379.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "'5.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.10%  "

$ws.Range("D17").Value = "'2.082.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Value = "'37.179.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").Value = "'16.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.83%  "

$ws.Range("D20").Value = "'75.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "

$ws.Range("D21").Value = "'0.0₃0936"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.50%  "

$ws.Range("D22").Value = "'5.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.97%  "

$ws.Range("D23").Value = "'239.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").Value = "  -2.34%  "

$ws.Range("E26").Value = "  +14.68%  "

$ws.Range("D27").Value = "'170.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("D28").Value = "'9.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("D29").Value = "'20.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("E30").Value = "  +3.25%  "

$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.83%  "

$ws.Range("D32").Value = "'4.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.43%  "

$ws.Range("D33").Value = "'0.0627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").Value = "'4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.48%  "

$ws.Range("D35").Value = "'0.0911"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  +27.35%  "

$ws.Range("D39").Value = "'1.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.09%  "

$ws.Range("D40").Value = "'1.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "'99.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").Value = "'4.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("D47").Value = "'4.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.35%  "

$ws.Range("E48").Value = "  +6.88%  "

$ws.Range("D49").Value = "'1.311.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").Value = "'2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("E51").Value = "  -0.39%  "

# Row 41/42: VeChain and InjectiveProtocol swapped order, with updated price/volume
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'18.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.17%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0228"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "

